$d = $word.ActiveDocument

$replacements = @(
    @("97×74=7178", "54×48=2592"),
    @("96×70=6720", "22×74=1628"),
    @("76×36=2736", "44×94=4136"),
    @("43×51=2193", "96×84=8064"),
    @("36×85=3060", "46×29=1334"),
    @("59×41=2419", "19×35=665"),
    @("81×42=3402", "41×12=492"),
    @("32×70=2240", "66×61=4026"),
    @("28×61=1708", "31×56=1736"),
    @("32×71=2272", "81×71=5751"),
    @("32×96=3072", "36×93=3348"),
    @("84×87=7308", "57×88=5016"),
    @("92×84=7728", "37×31=1147"),
    @("81×80=6480", "23×91=2093"),
    @("71×80=5680", "11×97=1067"),
    @("86×31=2666", "33×97=3201"),
    @("28×80=2240", "68×53=3604"),
    @("99×53=5247", "40×67=2680"),
    @("33×32=1056", "97×81=7857"),
    @("61×93=5673", "65×54=3510"),
    @("86×32=2752", "41×30=1230"),
    @("31×91=2821", "66×54=3564"),
    @("96×91=8736", "21×13=273"),
    @("75×15=1125", "82×54=4428"),
    @("31×88=2728", "77×97=7469")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
